$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell edits per diff ---
$ws.Range("Q62").Value = 0
$ws.Range("O1069").Value = 2
$ws.Range("R1071").Value = 0
$ws.Range("R1072").Value = 0

# --- Append new weekly rows 1073-1087 ---
# Row 1073
$ws.Range("A1073").Value = 45474
$ws.Range("A1073").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1073").Value = 506.7000122070312
$ws.Range("C1073").Value = 537.5
$ws.Range("D1073").Value = 502.8500061035156
$ws.Range("E1073").Value = 522.2999877929688
$ws.Range("F1073").Value = 517.3237915039062
$ws.Range("G1073").Value = 17182395
$ws.Range("H1073").Value = 2024
$ws.Range("I1073").Value = 7
$ws.Range("J1073").Value = 1
$ws.Range("K1073").Value = 0
$ws.Range("L1073").Value = 0
$ws.Range("M1073").Value = 0
$ws.Range("N1073").Value = 27
$ws.Range("O1073").Value = 0
$ws.Range("P1073").Value = 0
$ws.Range("Q1073").Value = 0

# Row 1074
$ws.Range("A1074").Value = 45481
$ws.Range("A1074").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1074").Value = 525.9500122070312
$ws.Range("C1074").Value = 540.4500122070312
$ws.Range("D1074").Value = 516.8499755859375
$ws.Range("E1074").Value = 525.5
$ws.Range("F1074").Value = 520.4933471679688
$ws.Range("G1074").Value = 10928015
$ws.Range("H1074").Value = 2024
$ws.Range("I1074").Value = 7
$ws.Range("J1074").Value = 8
$ws.Range("K1074").Value = 0
$ws.Range("L1074").Value = 0
$ws.Range("M1074").Value = 0
$ws.Range("N1074").Value = 28
$ws.Range("O1074").Value = 0
$ws.Range("P1074").Value = 0
$ws.Range("Q1074").Value = 0

# Row 1075
$ws.Range("A1075").Value = 45488
$ws.Range("A1075").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1075").Value = 527
$ws.Range("C1075").Value = 550.9000244140625
$ws.Range("D1075").Value = 521.5999755859375
$ws.Range("E1075").Value = 528.8499755859375
$ws.Range("F1075").Value = 523.8114013671875
$ws.Range("G1075").Value = 10626171
$ws.Range("H1075").Value = 2024
$ws.Range("I1075").Value = 7
$ws.Range("J1075").Value = 15
$ws.Range("K1075").Value = 0
$ws.Range("L1075").Value = 0
$ws.Range("M1075").Value = 0
$ws.Range("N1075").Value = 29
$ws.Range("O1075").Value = 0
$ws.Range("P1075").Value = 0
$ws.Range("Q1075").Value = 1

# Row 1076
$ws.Range("A1076").Value = 45495
$ws.Range("A1076").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1076").Value = 527.9500122070312
$ws.Range("C1076").Value = 553.4000244140625
$ws.Range("D1076").Value = 516.2999877929688
$ws.Range("E1076").Value = 540.5
$ws.Range("F1076").Value = 535.3504028320312
$ws.Range("G1076").Value = 11724170
$ws.Range("H1076").Value = 2024
$ws.Range("I1076").Value = 7
$ws.Range("J1076").Value = 22
$ws.Range("K1076").Value = 0
$ws.Range("L1076").Value = 0
$ws.Range("M1076").Value = 0
$ws.Range("N1076").Value = 30
$ws.Range("O1076").Value = 0
$ws.Range("P1076").Value = 0
$ws.Range("Q1076").Value = 0

# Row 1077
$ws.Range("A1077").Value = 45502
$ws.Range("A1077").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1077").Value = 545
$ws.Range("C1077").Value = 559.4000244140625
$ws.Range("D1077").Value = 531
$ws.Range("E1077").Value = 538.0999755859375
$ws.Range("F1077").Value = 532.9732666015625
$ws.Range("G1077").Value = 6952228
$ws.Range("H1077").Value = 2024
$ws.Range("I1077").Value = 7
$ws.Range("J1077").Value = 29
$ws.Range("K1077").Value = 0
$ws.Range("L1077").Value = 0
$ws.Range("M1077").Value = 0
$ws.Range("N1077").Value = 31
$ws.Range("O1077").Value = 0
$ws.Range("P1077").Value = 0
$ws.Range("Q1077").Value = 0

# Row 1078
$ws.Range("A1078").Value = 45509
$ws.Range("A1078").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1078").Value = 520
$ws.Range("C1078").Value = 549
$ws.Range("D1078").Value = 513.2999877929688
$ws.Range("E1078").Value = 537.4000244140625
$ws.Range("F1078").Value = 532.2799682617188
$ws.Range("G1078").Value = 9487489
$ws.Range("H1078").Value = 2024
$ws.Range("I1078").Value = 8
$ws.Range("J1078").Value = 5
$ws.Range("K1078").Value = 0
$ws.Range("L1078").Value = 0
$ws.Range("M1078").Value = 0
$ws.Range("N1078").Value = 32
$ws.Range("O1078").Value = 0
$ws.Range("P1078").Value = 0
$ws.Range("Q1078").Value = 0

# Row 1079
$ws.Range("A1079").Value = 45516
$ws.Range("A1079").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1079").Value = 537.0499877929688
$ws.Range("C1079").Value = 551
$ws.Range("D1079").Value = 533.5
$ws.Range("E1079").Value = 545.8499755859375
$ws.Range("F1079").Value = 540.6494140625
$ws.Range("G1079").Value = 3702684
$ws.Range("H1079").Value = 2024
$ws.Range("I1079").Value = 8
$ws.Range("J1079").Value = 12
$ws.Range("K1079").Value = 0
$ws.Range("L1079").Value = 0
$ws.Range("M1079").Value = 0
$ws.Range("N1079").Value = 33
$ws.Range("O1079").Value = 0
$ws.Range("P1079").Value = 0
$ws.Range("Q1079").Value = 0

# Row 1080
$ws.Range("A1080").Value = 45523
$ws.Range("A1080").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1080").Value = 545
$ws.Range("C1080").Value = 555
$ws.Range("D1080").Value = 522.75
$ws.Range("E1080").Value = 524.0499877929688
$ws.Range("F1080").Value = 519.05712890625
$ws.Range("G1080").Value = 6072164
$ws.Range("H1080").Value = 2024
$ws.Range("I1080").Value = 8
$ws.Range("J1080").Value = 19
$ws.Range("K1080").Value = 0
$ws.Range("L1080").Value = 0
$ws.Range("M1080").Value = 0
$ws.Range("N1080").Value = 34
$ws.Range("O1080").Value = 0
$ws.Range("P1080").Value = 0
$ws.Range("Q1080").Value = 0

# Row 1081
$ws.Range("A1081").Value = 45530
$ws.Range("A1081").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1081").Value = 528
$ws.Range("C1081").Value = 560.75
$ws.Range("D1081").Value = 519.9500122070312
$ws.Range("E1081").Value = 552.7999877929688
$ws.Range("F1081").Value = 547.533203125
$ws.Range("G1081").Value = 9409783
$ws.Range("H1081").Value = 2024
$ws.Range("I1081").Value = 8
$ws.Range("J1081").Value = 26
$ws.Range("K1081").Value = 0
$ws.Range("L1081").Value = 0
$ws.Range("M1081").Value = 0
$ws.Range("N1081").Value = 35
$ws.Range("O1081").Value = 0
$ws.Range("P1081").Value = 0
$ws.Range("Q1081").Value = 0

# Row 1082
$ws.Range("A1082").Value = 45537
$ws.Range("A1082").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1082").Value = 554.0999755859375
$ws.Range("C1082").Value = 570.3499755859375
$ws.Range("D1082").Value = 535.5499877929688
$ws.Range("E1082").Value = 542.3499755859375
$ws.Range("F1082").Value = 537.1827392578125
$ws.Range("G1082").Value = 14206758
$ws.Range("H1082").Value = 2024
$ws.Range("I1082").Value = 9
$ws.Range("J1082").Value = 2
$ws.Range("K1082").Value = 0
$ws.Range("L1082").Value = 0
$ws.Range("M1082").Value = 0
$ws.Range("N1082").Value = 36
$ws.Range("O1082").Value = 1
$ws.Range("P1082").Value = 0
$ws.Range("Q1082").Value = 0

# Row 1083
$ws.Range("A1083").Value = 45544
$ws.Range("A1083").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1083").Value = 538.5499877929688
$ws.Range("C1083").Value = 545
$ws.Range("D1083").Value = 516.0499877929688
$ws.Range("E1083").Value = 517.1500244140625
$ws.Range("F1083").Value = 512.222900390625
$ws.Range("G1083").Value = 4781949
$ws.Range("H1083").Value = 2024
$ws.Range("I1083").Value = 9
$ws.Range("J1083").Value = 9
$ws.Range("K1083").Value = 0
$ws.Range("L1083").Value = 0
$ws.Range("M1083").Value = 0
$ws.Range("N1083").Value = 37
$ws.Range("O1083").Value = 0
$ws.Range("P1083").Value = 0
$ws.Range("Q1083").Value = 0

# Row 1084
$ws.Range("A1084").Value = 45551
$ws.Range("A1084").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1084").Value = 518
$ws.Range("C1084").Value = 562.75
$ws.Range("D1084").Value = 517.1500244140625
$ws.Range("E1084").Value = 538.9000244140625
$ws.Range("F1084").Value = 538.9000244140625
$ws.Range("G1084").Value = 27970619
$ws.Range("H1084").Value = 2024
$ws.Range("I1084").Value = 9
$ws.Range("J1084").Value = 16
$ws.Range("K1084").Value = 0
$ws.Range("L1084").Value = 0
$ws.Range("M1084").Value = 0
$ws.Range("N1084").Value = 38
$ws.Range("O1084").Value = 0
$ws.Range("P1084").Value = 0
$ws.Range("Q1084").Value = 0

# Row 1085
$ws.Range("A1085").Value = 45558
$ws.Range("A1085").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1085").Value = 538.5499877929688
$ws.Range("C1085").Value = 558.5
$ws.Range("D1085").Value = 534
$ws.Range("E1085").Value = 549.5
$ws.Range("F1085").Value = 549.5
$ws.Range("G1085").Value = 8698991
$ws.Range("H1085").Value = 2024
$ws.Range("I1085").Value = 9
$ws.Range("J1085").Value = 23
$ws.Range("K1085").Value = 0
$ws.Range("L1085").Value = 0
$ws.Range("M1085").Value = 0
$ws.Range("N1085").Value = 39
$ws.Range("O1085").Value = 0
$ws.Range("P1085").Value = 0
$ws.Range("Q1085").Value = 0

# Row 1086
$ws.Range("A1086").Value = 45565
$ws.Range("A1086").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1086").Value = 549
$ws.Range("C1086").Value = 569.5499877929688
$ws.Range("D1086").Value = 543
$ws.Range("E1086").Value = 549.3499755859375
$ws.Range("F1086").Value = 549.3499755859375
$ws.Range("G1086").Value = 9789551
$ws.Range("H1086").Value = 2024
$ws.Range("I1086").Value = 9
$ws.Range("J1086").Value = 30
$ws.Range("K1086").Value = 0
$ws.Range("L1086").Value = 0
$ws.Range("M1086").Value = 0
$ws.Range("N1086").Value = 40
$ws.Range("O1086").Value = 0
$ws.Range("P1086").Value = 0
$ws.Range("Q1086").Value = 0

# Row 1087
$ws.Range("A1087").Value = 45572
$ws.Range("A1087").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1087").Value = 553
$ws.Range("C1087").Value = 553.0499877929688
$ws.Range("D1087").Value = 529.9000244140625
$ws.Range("E1087").Value = 540.5499877929688
$ws.Range("F1087").Value = 540.5499877929688
$ws.Range("G1087").Value = 10660063
$ws.Range("H1087").Value = 2024
$ws.Range("I1087").Value = 10
$ws.Range("J1087").Value = 7
$ws.Range("K1087").Value = 0
$ws.Range("L1087").Value = 0
$ws.Range("M1087").Value = 0
$ws.Range("N1087").Value = 41
$ws.Range("O1087").Value = 0
$ws.Range("P1087").Value = 0
$ws.Range("Q1087").Value = 0
